$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

# row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 143.06667
$ws.Range("I2").Value = 146.64285
$ws.Range("K2").Value = 146.64285
$ws.Range("M2").Value = -33.64285000000001

# row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 471.3889
$ws.Range("I41").Value = 301.5
$ws.Range("K41").Value = 301.5
$ws.Range("M41").Value = 138.5

# row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 104.111115
$ws.Range("I55").Value = 57
$ws.Range("J55").Value = 127.666664
$ws.Range("K55").Value = 57
$ws.Range("L55").Value = 127.666664
$ws.Range("M55").Value = 157
$ws.Range("N55").Value = -555.666664

# row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 866.7619
$ws.Range("I111").Value = 896.64703
$ws.Range("J111").Value = 739.75
$ws.Range("K111").Value = 2689.94109
$ws.Range("L111").Value = 2219.25
$ws.Range("M111").Value = 377.0589100000002
$ws.Range("N111").Value = -8353.25

# row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 2500
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 8181.4546
$ws.Range("I116").Value = 5843.154
$ws.Range("J116").Value = 11559
$ws.Range("K116").Value = 5843.154
$ws.Range("L116").Value = 11559
$ws.Range("M116").Value = -2401.154
$ws.Range("N116").Value = -18443

# row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 858.5714
$ws.Range("I118").Value = 858.5714
$ws.Range("K118").Value = 2575.7142
$ws.Range("M118").Value = -918.7142000000003

# row 133 (Leve Item ID 41856)
$ws.Range("H133").Value = 52280
$ws.Range("J133").Value = 52280
$ws.Range("L133").Value = 52280
$ws.Range("N133").Value = -62400

# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3772.8572
$ws.Range("I138").Value = 2931.5
$ws.Range("J138").Value = 4481.3687
$ws.Range("K138").Value = 8794.5
$ws.Range("L138").Value = 13444.1061
$ws.Range("M138").Value = -3654.5
$ws.Range("N138").Value = -23724.1061


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

# row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1730.4722
$ws.Range("I2").Value = 1694
$ws.Range("J2").Value = 1912.8334
$ws.Range("K2").Value = 1694
$ws.Range("L2").Value = 1912.8334
$ws.Range("M2").Value = -1581
$ws.Range("N2").Value = -2138.8334

# row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 2125.423
$ws.Range("J110").Value = 2084.8333
$ws.Range("L110").Value = 2084.8333
$ws.Range("N110").Value = -6174.8333

# row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1730.4722
$ws.Range("I116").Value = 1694
$ws.Range("J116").Value = 1912.8334
$ws.Range("K116").Value = 1694
$ws.Range("L116").Value = 1912.8334
$ws.Range("M116").Value = 600
$ws.Range("N116").Value = -6500.8334

# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3615.5535
$ws.Range("I132").Value = 3615.5535
$ws.Range("K132").Value = 10846.6605
$ws.Range("M132").Value = -8316.6605


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")

# row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1730.4722
$ws.Range("I3").Value = 1694
$ws.Range("J3").Value = 1912.8334
$ws.Range("K3").Value = 1694
$ws.Range("L3").Value = 1912.8334
$ws.Range("M3").Value = -1580
$ws.Range("N3").Value = -2140.8334

# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 558451.4
$ws.Range("I86").Value = 911503.9399999999
$ws.Range("K86").Value = 911503.9399999999
$ws.Range("M86").Value = -910380.9399999999

# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 558451.4
$ws.Range("I89").Value = 911503.9399999999
$ws.Range("K89").Value = 4557519.699999999
$ws.Range("M89").Value = -4551903.699999999

# row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 1672
$ws.Range("I107").Value = 1672
$ws.Range("K107").Value = 1672
$ws.Range("M107").Value = 248

# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4500.3228
$ws.Range("I134").Value = 4250.3667
$ws.Range("K134").Value = 12751.1001
$ws.Range("M134").Value = -10216.1001


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3034.721
$ws.Range("I31").Value = 1747.0714
$ws.Range("J31").Value = 5438.3335
$ws.Range("K31").Value = 1747.0714
$ws.Range("L31").Value = 5438.3335
$ws.Range("M31").Value = -1452.0714
$ws.Range("N31").Value = -6028.3335

# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3034.721
$ws.Range("I34").Value = 1747.0714
$ws.Range("J34").Value = 5438.3335
$ws.Range("K34").Value = 1747.0714
$ws.Range("L34").Value = 5438.3335
$ws.Range("M34").Value = -1545.0714
$ws.Range("N34").Value = -5842.3335

# row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 2006
$ws.Range("I58").Value = 1963
$ws.Range("K58").Value = 1963
$ws.Range("M58").Value = -1760

# row 95 (Leve Item ID 18192)
$ws.Range("H95").Value = 28110.75
$ws.Range("J95").Value = 28110.75
$ws.Range("L95").Value = 28110.75
$ws.Range("N95").Value = -33602.75

# row 96 (Leve Item ID 18193)
$ws.Range("H96").Value = 21724
$ws.Range("J96").Value = 21724
$ws.Range("L96").Value = 21724
$ws.Range("N96").Value = -27216

# row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 10831.333
$ws.Range("I99").Value = 10597.6
$ws.Range("J99").Value = 12000
$ws.Range("K99").Value = 10597.6
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = -9099.6
$ws.Range("N99").Value = -14996

# row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 10831.333
$ws.Range("I126").Value = 10597.6
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 31792.8
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -29322.8
$ws.Range("N126").Value = -40940

# row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 5973.815
$ws.Range("I132").Value = 4692
$ws.Range("K132").Value = 14076
$ws.Range("M132").Value = -11546

# row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 2006
$ws.Range("I136").Value = 1963
$ws.Range("K136").Value = 5889
$ws.Range("M136").Value = -3339


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")

# row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 52976190
$ws.Range("I4").Value = 55451250
$ws.Range("K4").Value = 166353750
$ws.Range("M4").Value = -166353638

# row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 2082.7058
$ws.Range("I5").Value = 921.7
$ws.Range("J5").Value = 2566.4583
$ws.Range("K5").Value = 2765.1
$ws.Range("L5").Value = 7699.374899999999
$ws.Range("M5").Value = -2653.1
$ws.Range("N5").Value = -7923.374899999999

# row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 2082.7058
$ws.Range("I135").Value = 921.7
$ws.Range("J135").Value = 2566.4583
$ws.Range("K135").Value = 8295.300000000001
$ws.Range("L135").Value = 23098.1247
$ws.Range("M135").Value = -5760.300000000001
$ws.Range("N135").Value = -28168.1247


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

# row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2466.423
$ws.Range("I132").Value = 2406.7
$ws.Range("J132").Value = 2665.5
$ws.Range("K132").Value = 7220.099999999999
$ws.Range("L132").Value = 7996.5
$ws.Range("M132").Value = -4690.099999999999
$ws.Range("N132").Value = -13056.5


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2528.6
$ws.Range("J22").Value = 2239.4285
$ws.Range("L22").Value = 2239.4285
$ws.Range("N22").Value = -2829.4285

# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2528.6
$ws.Range("J27").Value = 2239.4285
$ws.Range("L27").Value = 2239.4285
$ws.Range("N27").Value = -2453.4285

# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1418.6666
$ws.Range("I46").Value = 1100
$ws.Range("J46").Value = 1518.25
$ws.Range("K46").Value = 1100
$ws.Range("L46").Value = 1518.25
$ws.Range("M46").Value = -912
$ws.Range("N46").Value = -1894.25

# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3328.9167
$ws.Range("I132").Value = 2994.3
$ws.Range("K132").Value = 8982.900000000001
$ws.Range("M132").Value = -6452.900000000001


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

# row 34 (Leve Item ID 3349)
$ws.Range("H34").Value = 1700
$ws.Range("I34").Value = 1700
$ws.Range("K34").Value = 1700
$ws.Range("M34").Value = -1497

# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 6301.3335
$ws.Range("I132").Value = 6606.615
$ws.Range("K132").Value = 19819.845
$ws.Range("M132").Value = -17289.845

# row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 687.64703
$ws.Range("I136").Value = 687.64703
$ws.Range("K136").Value = 2062.94109
$ws.Range("M136").Value = 487.0589100000002

